# Generate Report for Handback
# Refresh the "Latest HO Xliff Generate Date" / handoff & handback datetime
# stamps for the afddc8bf-b5fd-4493-8abe-fb20b78e3889 row across the
# Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Sheets.Item("Overview")
$zhcn     = $wb.Sheets.Item("zh-cn")
$dede     = $wb.Sheets.Item("de-de")

# Overview: Latest HO Xliff Generate Date
$overview.Range("G4").Value = "2017-02-21 03:56:54"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime
$zhcn.Range("H4").Value = "2017-02-21 03:56:37"
$zhcn.Range("L4").Value = "2017-02-21 03:57:33"

# de-de: Correspond Handoff Datetime (mirrors the Overview value) /
# Correspond Handback DateTime
$dede.Range("H4").Value = "2017-02-21 03:56:54"
$dede.Range("L4").Value = "2017-02-21 03:57:55"
